$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44315
$ws.Range("I3").Value = "Especial"
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = 30000
$ws.Range("N3").Value = "$/caja 20 kilos empedrada"
$ws.Range("P3").Value = 1500
$ws.Range("Q3").Value = 20

# Row 4
$ws.Range("D4").Value = 44315
$ws.Range("N4").Value = "$/caja 15 kilos granel"

# Row 5
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("N5").Value = "$/caja 15 kilos empedrada"
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = 15

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 20

# Row 7
$ws.Range("D7").Value = 44280
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 25000
$ws.Range("N7").Value = "$/caja 18 kilos empedrada"
$ws.Range("P7").Value = 1389
$ws.Range("Q7").Value = 18
